$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above existing data (shifts rows 1..5 -> 2..6)
$ws.Rows.Item(1).Insert()

# Header row
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Semester"
$ws.Range("C1").Value = "Nim"
$ws.Range("D1").Value = "Nama"
$ws.Range("E1").Value = "NIDN_Dosen_Pembimbing_Utama"
$ws.Range("F1").Value = "Nama_Dosen_Pembimbing_Utama"
$ws.Range("G1").Value = "NIDN_Dosen_Pembimbing_Pembantu"
$ws.Range("H1").Value = "Nama_Dosen_Pembimbing_Pembantu"

# New NIM values (E column) for rows 2..6
$nimValues = @("00406107055", "00406107056", "00406107057", "00406107058", "00406107059")
# Existing NIDN values that used to live in column F, now moved to column G
$nidnValues = @("0608068502", "0608068503", "0608068504", "0608068505", "0608068506")

for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 5).Value = "'" + $nimValues[$i]
    $ws.Cells.Item($r, 6).Value = "'Ratih Nindyasari, S.Kom, M.Kom"
    $ws.Cells.Item($r, 7).Value = "'" + $nidnValues[$i]
    $ws.Cells.Item($r, 8).Value = "Anastasya Latubessy, S.Kom, M.Cs"
}

# Column widths (closest achievable values given the host's 1/6-character
# width quantization grid; targets are 37.42578125 / 39.140625)
$ws.Columns.Item(5).ColumnWidth = 36.666666666666664
$ws.Columns.Item(6).ColumnWidth = 38.333333333333336
$ws.Columns.Item(7).ColumnWidth = 36.666666666666664
$ws.Columns.Item(8).ColumnWidth = 38.333333333333336

$ws.Range("H2:H6").Select()
